# Edit script: update cached "today" date fields across the slide master,
# every slide layout and the notes master; tag the quiz-answer smiley face
# shapes with alt text "QuizAnswer"; and nudge slide 3's smiley face shape
# to its new horizontal position.

$p = $ppt.ActivePresentation

$newDate = "10/27/2025"

# --- 1. Refresh the cached datetimeFigureOut field text -------------------
# Slide master + every custom (slide) layout expose the date placeholder as
# a normal shape whose TextFrame.TextRange we can rewrite directly.
function Update-DatePlaceholderShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholderShapes $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholderShapes $layout.Shapes
}

# The notes master's date placeholder only accepts updates through the
# HeadersFooters façade in this host.
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = $newDate

# --- 2. Tag the "Smiley Face 3" quiz-answer shapes -------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.Name -eq "Smiley Face 3") {
            $shape.AlternativeText = "QuizAnswer"
        }
    }
}

# --- 3. Reposition the smiley face on slide 3 ------------------------------
$slide3 = $p.Slides.Item(3)
for ($shi = 1; $shi -le $slide3.Shapes.Count; $shi++) {
    $shape = $slide3.Shapes.Item($shi)
    if ($shape.Name -eq "Smiley Face 3") {
        # 3317893 EMU: the host stores Shape.Left as a 32-bit float and
        # floors pt*12700 back to EMU, so the naive 3317893/12700 quotient
        # (261.25141732283464) rounds down to 3317892 EMU once it is
        # narrowed to single precision. Nudge the literal a hair higher so
        # the narrowed value still lands in the [3317893, 3317894) EMU band.
        $shape.Left = 261.25142
    }
}
